$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Student Summary": add "Course Code:" / "DSPC602" and "Max marks" / 40
# rows (new rows 11 & 12, between the Staff Name row and the Attribute/Value
# table), update the marks-distribution labels and fix the rounded average.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Student Summary")

# Clone formatting of the existing "Staff Name :" row onto the two new rows
# so the new cells pick up style index 2 (same as the rest of the info block).
$wsSummary.Range("A10:C10").Copy($wsSummary.Range("A11"))
$wsSummary.Range("A10:C10").Copy($wsSummary.Range("A12"))

$wsSummary.Range("A11").Value = ""
$wsSummary.Range("B11").Value = "Course Code:"
$wsSummary.Range("C11").Value = "DSPC602"

$wsSummary.Range("A12").Value = ""
$wsSummary.Range("B12").Value = "Max marks"
$wsSummary.Range("C12").Value = 40

# Relabel the marks-distribution buckets (rows 17-20 keep their row numbers).
$wsSummary.Range("A17").Value = "Average Marks"
$wsSummary.Range("B17").Value = 19.88
$wsSummary.Range("A18").Value = "Less Than 40%"
$wsSummary.Range("A19").Value = "Between 40 % - 75 %"
$wsSummary.Range("A20").Value = "More than 75%"

# ---------------------------------------------------------------------------
# Sheet "Fast Learners": fix the roll-number/name mismatch for the bottom
# block (rows 9-15) and drop the trailing three rows that belonged in the
# absentee list, not here.
# ---------------------------------------------------------------------------
$wsFast = $wb.Worksheets.Item("Fast Learners")

$oldRows = @{}
for ($r = 9; $r -le 15; $r++) {
    $oldRows[$r] = @(
        $wsFast.Cells.Item($r, 1).Value2,
        $wsFast.Cells.Item($r, 2).Value2,
        $wsFast.Cells.Item($r, 3).Value2,
        $wsFast.Cells.Item($r, 4).Value2
    )
}

# Order (by former row number) that the corrected roll-no/name pairs should
# appear in, starting again at row 9.
$newOrder = @(13, 15, 14, 12, 9, 11, 10)
$destRow = 9
foreach ($srcRow in $newOrder) {
    $vals = $oldRows[$srcRow]
    $wsFast.Cells.Item($destRow, 1).Value = $vals[0]
    $wsFast.Cells.Item($destRow, 2).Value = $vals[1]
    $wsFast.Cells.Item($destRow, 3).Value = $vals[2]
    $wsFast.Cells.Item($destRow, 4).Value = $vals[3]
    $destRow = $destRow + 1
}

# Remove the three rows that no longer belong on this sheet (former rows
# 16-18 -- DEEPAKRAGAVAN J / PREETHIGA S / KRISHNAKUMAR S).
$wsFast.Rows.Item(16).Delete()
$wsFast.Rows.Item(16).Delete()
$wsFast.Rows.Item(16).Delete()
